{"js": "// Delete the whole list paragraph \"Take out uniquely renumber\".\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text.trim() === \"Take out uniquely renumber\") {\n    p.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Delete the whole list paragraph \"Take out uniquely renumber\".\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\nfor ($i = $count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.Trim() -eq \"Take out uniquely renumber\") {\n        $p.Range.Delete()\n    }\n}\n"}
